$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 146, shifting existing rows 146..252 down to 147..253.
$ws.Rows(146).Insert()

# Populate the newly inserted row 146 with the new weekly data record.
$ws.Cells.Item(146, 1).Value = 11
$ws.Cells.Item(146, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(146, 3).Value = "Bíobío"
$ws.Cells.Item(146, 4).Value = 44978
$ws.Cells.Item(146, 5).Value = 8
$ws.Cells.Item(146, 6).Value = 100112003
$ws.Cells.Item(146, 7).Value = "Ajo"
$ws.Cells.Item(146, 8).Value = "Chino"
$ws.Cells.Item(146, 9).Value = "Primera"
$ws.Cells.Item(146, 10).Value = 400
$ws.Cells.Item(146, 11).Value = 15000
$ws.Cells.Item(146, 12).Value = 16000
$ws.Cells.Item(146, 13).Value = 15500
$ws.Cells.Item(146, 14).Value = '$/caja 10 kilos'
$ws.Cells.Item(146, 15).Value = "China"
$ws.Cells.Item(146, 16).Value = 1550
$ws.Cells.Item(146, 17).Value = 10
$ws.Cells.Item(146, 18).Value = "Hortaliza"
